$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append two new daily rows (2025-10-18) for both stations, matching the
# formatting of the preceding data rows (copy formats only, so the existing
# style indices are reused rather than minting new ones).
$ws.Range("A35").Copy()
$ws.Range("A36").PasteSpecial(-4122)
$ws.Range("A37").PasteSpecial(-4122)

$ws.Range("C35:F35").Copy()
$ws.Range("C36:F36").PasteSpecial(-4122)
$ws.Range("C37:F37").PasteSpecial(-4122)

$ws.Cells.Item(36, 1).Value = 45948
$ws.Cells.Item(36, 2).Value = "四方坪站"
$ws.Cells.Item(36, 3).Value = 10067.02
$ws.Cells.Item(36, 4).Value = 8291.7999999999993
$ws.Cells.Item(36, 5).Value = 3486
$ws.Cells.Item(36, 6).Value = 434

$ws.Cells.Item(37, 1).Value = 45948
$ws.Cells.Item(37, 2).Value = "高岭站"
$ws.Cells.Item(37, 3).Value = 3944.96
$ws.Cells.Item(37, 4).Value = 3210.33
$ws.Cells.Item(37, 5).Value = 989.79
$ws.Cells.Item(37, 6).Value = 149

# Update the active selection, matching the recorded view state.
$ws.Range("L31").Select()
